$wb = $excel.ActiveWorkbook

# New row (row 22) data for each worksheet, in workbook sheet order:
#   1 = ROW35-FE-LIFTER
#   2 = ROW35-MID-LIFTER
#   3 = ROW02-FE-LIFTER
#   4 = ROW02-MID-LIFTER

$rows = @(
    @{
        Sheet = 1
        A = 45735.27706204861
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x82"
        E = "0xd"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 386
        I = 13
    },
    @{
        Sheet = 2
        A = 45735.12711385416
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x86"
        E = "0xe"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 390
        I = 14
    },
    @{
        Sheet = 3
        A = 45735.27407221065
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x01,0x82"
        E = "0x3"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 386
        I = 3
    },
    @{
        Sheet = 4
        A = 45735.33468563658
        B = "0x01,0x90"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x01,0x86"
        E = "0x3"
        F = 400
        G = [double]"9.85046333984776e+23"
        H = 390
        I = 3
    }
)

foreach ($row in $rows) {
    $ws = $wb.Worksheets.Item($row.Sheet)

    # Mirror the date-time number format used by column A of the other rows.
    $ws.Range("A22").NumberFormat = $ws.Range("A21").NumberFormat

    $ws.Range("A22").Value = $row.A
    $ws.Range("B22").Value = $row.B
    $ws.Range("C22").Value = $row.C
    $ws.Range("D22").Value = $row.D
    $ws.Range("E22").Value = $row.E
    $ws.Range("F22").Value = $row.F
    $ws.Range("G22").Value = $row.G
    $ws.Range("H22").Value = $row.H
    $ws.Range("I22").Value = $row.I
}
